# Adicionado Scaller para garantir que valores fiquem entre min e max.
# Extends the DOE Full Factorial table by adding a new factor "Temperatura"
# (column F) with two levels (20 and 50), and adjusts the existing
# "Discreto 2" (column E) values so the design becomes a full factorial
# across all factors (A1:E17 -> A1:F33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F
$ws.Range("F1").Value = "Temperatura"

# Copy the header style (bold + border) from an existing header cell onto
# the new header cell.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# The 8 base combinations of Viscosity (B), Densidade (C) and Numero de
# Pratos (D) that repeat for every level of E (Discreto 2) / F (Temperatura)
$combos = @(
    @(717.7573186524264, 1000, 10),
    @(882.2426813475736, 1000, 10),
    @(717.7573186524264, 2000, 10),
    @(882.2426813475736, 2000, 10),
    @(717.7573186524264, 1000, 15),
    @(882.2426813475736, 1000, 15),
    @(717.7573186524264, 2000, 15),
    @(882.2426813475736, 2000, 15)
)

# Blocks describe how the sheet is laid out: starting row, the value that
# "Discreto 2" (E) should take for that block and the value "Temperatura"
# (F) should take for that block.
$blocks = @(
    @{ StartRow = 2;  E = 2;  F = 20 },
    @{ StartRow = 10; E = 10; F = 20 },
    @{ StartRow = 18; E = 2;  F = 50 },
    @{ StartRow = 26; E = 10; F = 50 }
)

foreach ($block in $blocks) {
    $startRow = $block.StartRow
    $eVal = $block.E
    $fVal = $block.F

    for ($i = 0; $i -lt $combos.Length; $i++) {
        $row = $startRow + $i
        $combo = $combos[$i]

        $ws.Cells.Item($row, 1).Value = $row - 1
        $ws.Cells.Item($row, 2).Value = $combo[0]
        $ws.Cells.Item($row, 3).Value = $combo[1]
        $ws.Cells.Item($row, 4).Value = $combo[2]
        $ws.Cells.Item($row, 5).Value = $eVal
        $ws.Cells.Item($row, 6).Value = $fVal
    }
}

# Apply the same style used by column A's data cells (bold + border) to the
# "Simulation" cells of the newly added rows (18-33).
$ws.Range("A2").Copy()
$ws.Range("A18:A33").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
